# Electricity fixes for solid waste and geothermal MPCbS and SYC
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsMPCbS = $wb.Worksheets.Item("MPCbS")

# Geothermal max potential capacity (MW) -> 0
$wsData.Range("B10").Value = 0

# Municipal solid waste max potential capacity (MW) -> 0
$wsData.Range("B17").Value = 0

# Make MPCbS the active / selected sheet (was About)
$wsMPCbS.Activate()
$wsMPCbS.Select()

$wb.Save()
